$p = $ppt.ActivePresentation

# --- Slide 1: bump the date in the "Alison Smith | <date>" subtitle. ---
$s1 = $p.Slides.Item(1)
$tr1 = $s1.Shapes.Item(5).TextFrame.TextRange
$old1 = "Alison Smith | November 17, 2025"
$new1 = "Alison Smith | November 18, 2025"
$tr1.Characters(1, $old1.Length).Text = $new1

# --- Slide 5: drop the markdown-style asterisks around each phase's ---
# --- month range (e.g. " *(Months 1-2)*" -> " (Months 1-2)").       ---
$s5 = $p.Slides.Item(5)
$tr5 = $s5.Shapes.Item(3).TextFrame.TextRange

$replacements = @(
    @{ Old = " *(Months 1-2)*"; New = " (Months 1-2)" },
    @{ Old = " *(Months 3-4)*"; New = " (Months 3-4)" },
    @{ Old = " *(Months 5-6)*"; New = " (Months 5-6)" }
)

# Re-read the shape's full text before each replacement so that the
# character offsets always reflect the text as it currently stands
# (earlier replacements shrink the string and shift later offsets).
foreach ($rep in $replacements) {
    $full5 = $tr5.Text
    $idx = $full5.IndexOf($rep.Old)
    if ($idx -lt 0) {
        throw "Could not find '$($rep.Old)' in shape text"
    }
    $tr5.Characters($idx + 1, $rep.Old.Length).Text = $rep.New
}
